# Adding updated gym workout data (January 2018, week 1, Saturday workout)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All eight new rows share the same date / week / month / year / day values
$workoutDate = Get-Date -Year 2018 -Month 1 -Day 6 -Hour 0 -Minute 0 -Second 0

$newExercises = @(
    @{ Name = "Bench Press";     Weight = 95;  Sets = 5; Reps = 5;  Area = "Chest" },
    @{ Name = "Barbell Row";     Weight = 85;  Sets = 5; Reps = 5;  Area = "Back" },
    @{ Name = "Shoulder Press";  Weight = 26;  Sets = 4; Reps = 8;  Area = "Shoulders" },
    @{ Name = "Shoulder Shrug";  Weight = 26;  Sets = 4; Reps = 8;  Area = "Shoulders" },
    @{ Name = "Rear Delt";       Weight = 86;  Sets = 4; Reps = 8;  Area = "Back" },
    @{ Name = "Sit ups";         Weight = 0;   Sets = 5; Reps = 12; Area = "Core" },
    @{ Name = "Russian Twists";  Weight = 10;  Sets = 4; Reps = 12; Area = "Core" },
    @{ Name = "Pull-Ups";        Weight = 105; Sets = 5; Reps = 5;  Area = "Back" }
)

$startRow = 314
$exerciseId = 313
$dateId = 37

for ($i = 0; $i -lt $newExercises.Count; $i++) {
    $r = $startRow + $i
    $ex = $newExercises[$i]

    $ws.Cells.Item($r, 1).Value = $exerciseId      # A: ExerciseId
    $ws.Cells.Item($r, 2).Value = $dateId           # B: DateId
    $ws.Cells.Item($r, 3).Value = $workoutDate      # C: Exercise Date
    $ws.Cells.Item($r, 4).Value = 1                 # D: Exercise Week
    $ws.Cells.Item($r, 5).Value = "January"         # E: Exercise Month
    $ws.Cells.Item($r, 6).Value = 2018              # F: Exercise Year
    $ws.Cells.Item($r, 7).Value = "Saturday"        # G: Exercise Day
    $ws.Cells.Item($r, 8).Value = $ex.Name          # H: Exercise Name
    $ws.Cells.Item($r, 9).Value = $ex.Weight        # I: Weight
    $ws.Cells.Item($r, 10).Value = $ex.Sets         # J: Sets
    $ws.Cells.Item($r, 11).Value = $ex.Reps         # K: Reps
    $ws.Cells.Item($r, 12).Value = $ex.Area         # L: TrainingArea

    $exerciseId = $exerciseId + 1
}

# Match the author's final selection/scroll position recorded in the sheet view
$ws.Range("C322").Select()
